$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-08-31"
$ws.Range("I1").Value = "2022 (through 08-31)"
$ws.Range("I9").Value = 167
$ws.Range("I14").Value = 1138
